# Scenarios.xlsx - "Setting SlipRatio to 1" commit
#
# SlipRatio<1 is only required in a 2D situation; for the paper we simplify
# by setting SlipRatio (column J) to 1 for (almost) every scenario row, and
# flip the corresponding "Run" flag (column B) to TRUE so those scenarios are
# included in the run. Row 10 is the exception - it keeps the old 0.25 value
# that every other row used to have (its Run flag was already TRUE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Run flags (column B): turn every scenario on ------------------------
$ws.Range("B3:B50").Value = $true

# --- SlipRatio (column J): set to 1 for every row except row 10 ----------
$ws.Range("J3:J50").Value = 1
$ws.Range("J10").Value = 0.25

# --- Restore the last-used selection recorded in the sheet view ----------
$ws.Range("S10").Select() | Out-Null
